$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Create the new "centered / middle" data style on the anchor cell A2 and
#    set its value (also starts off column A). Doing both alignment
#    properties on this one cell creates exactly one new cellXfs entry.
# ---------------------------------------------------------------------------
$anchor = $ws.Range("A2")
$anchor.Value = "Oneway"
$anchor.HorizontalAlignment = -4108
$anchor.VerticalAlignment = -4108

# Propagate that style to the whole data block via copy / paste-special so we
# do not generate additional style entries.
$anchor.Copy()
$ws.Range("A2:P8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the rest of column A (flight type) for rows 3-7.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Return"
$ws.Range("A4").Value = "Multi"
$ws.Range("A5").Value = "Oneway"
$ws.Range("A6").Value = "Return"
$ws.Range("A7").Value = "Multi"

# ---------------------------------------------------------------------------
# 3) From / Destination columns (C & D) for rows 2-7.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Melbourne"
$ws.Range("D2").Value = "Sydney"

$ws.Range("C3").Value = "Melbourne"
$ws.Range("D3").Value = "Tasmania"

$ws.Range("C4").Value = "Melbourne"
$ws.Range("D4").Value = "Adelaide"

$ws.Range("C5").Value = "Melbourne"
$ws.Range("D5").Value = "Adelaide"

$ws.Range("C6").Value = "Sydney"
$ws.Range("D6").Value = "Melbourne"

$ws.Range("C7").Value = "Sydney"
$ws.Range("D7").Value = "Hobart"

# ---------------------------------------------------------------------------
# 4) Header: F1 "Arrival Date" -> "Return Date"
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Return Date"

# ---------------------------------------------------------------------------
# 5) Class column (J) for rows 2-7.
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "Economy"
$ws.Range("J3").Value = "First Class"
$ws.Range("J4").Value = "Business"
$ws.Range("J5").Value = "Premium Economy"
$ws.Range("J6").Value = "Economy"
$ws.Range("J7").Value = "Economy"

# ---------------------------------------------------------------------------
# 6) Remaining numeric / date / reused-text columns.
# ---------------------------------------------------------------------------
# Row 2
$ws.Range("E2").Value = 25012022
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

# Row 3
$ws.Range("E3").Value = 25022022
$ws.Range("F3").Value = 27022022
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4
$ws.Range("B4").Value = 2
$ws.Range("E4").Value = 25012022
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = "Sydney"
$ws.Range("L4").Value = 27012022
$ws.Range("M4").Value = "Melbourne"
$ws.Range("N4").Value = 28012022
# F4 has no cell at all in the final sheet.
$ws.Range("F4").Clear()

# Row 5
$ws.Range("E5").Value = 25012022
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1

# Row 6
$ws.Range("E6").Value = 25022022
$ws.Range("F6").Value = 27022022
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 2

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("E7").Value = 25012022
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("K7").Value = "Melbourne"
$ws.Range("L7").Value = 27012022

# ---------------------------------------------------------------------------
# 7) Column J width.
# ---------------------------------------------------------------------------
$ws.Columns("J").ColumnWidth = 18.14

# ---------------------------------------------------------------------------
# 8) Selection.
# ---------------------------------------------------------------------------
$ws.Range("A2:P7").Select()
